# Update evaluation metrics for versions 3.13, 3.14 and 3.15 on the
# "Eval (binary)" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eval (binary)")

# Row 255 - 3.13 (Random Forest)
$ws.Range("B255").Value = 0.88
$ws.Range("D255").Value = 0.88
$ws.Range("F255").Value = 0.88
$ws.Range("G255").Value = 0.88

# Row 256 - 3.13 (XGBoost)
$ws.Range("B256").Value = 0.86
$ws.Range("D256").Value = 0.89
$ws.Range("E256").Value = 0.88
$ws.Range("F256").Value = 0.94
$ws.Range("G256").Value = 0.91

# Row 257 - 3.13 (LightGBM)
$ws.Range("B257").Value = 0.86
$ws.Range("C257").Value = 0.93
$ws.Range("D257").Value = 0.89
$ws.Range("F257").Value = 0.95

# Row 259 - 3.14 (Random Forest)
$ws.Range("B259").Value = 0.89
$ws.Range("C259").Value = 0.89
$ws.Range("D259").Value = 0.89
$ws.Range("E259").Value = 0.9
$ws.Range("F259").Value = 0.92
$ws.Range("G259").Value = 0.91

# Row 260 - 3.14 (XGBoost)
$ws.Range("B260").Value = 0.87
$ws.Range("C260").Value = 0.93
$ws.Range("E260").Value = 0.88
$ws.Range("G260").Value = 0.91

# Row 261 - 3.14 (LightGBM)
$ws.Range("C261").Value = 0.95
$ws.Range("E261").Value = 0.88
$ws.Range("F261").Value = 0.95
$ws.Range("G261").Value = 0.91

# Row 263 - 3.15 (Random Forest)
$ws.Range("B263").Value = 0.88
$ws.Range("D263").Value = 0.88
$ws.Range("E263").Value = 0.89
$ws.Range("F263").Value = 0.9

# Row 264 - 3.15 (XGBoost)
$ws.Range("E264").Value = 0.87
$ws.Range("F264").Value = 0.93

# Row 265 - 3.15 (LightGBM)
$ws.Range("B265").Value = 0.86
$ws.Range("D265").Value = 0.89
$ws.Range("E265").Value = 0.85

# Update the selected cell on the sheet view.
$ws.Range("G266").Select()
